$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 481513.72
$ws.Range("J17").Value = 503589.4
$ws.Range("L17").Value = 1510768.2
$ws.Range("N17").Value = -1511104.2

$ws.Range("H53").Value = 47988.81
$ws.Range("I53").Value = 407.7
$ws.Range("J53").Value = 91244.37
$ws.Range("K53").Value = 407.7
$ws.Range("L53").Value = 91244.37
$ws.Range("M53").Value = 229.3
$ws.Range("N53").Value = -92518.37

$ws.Range("H106").Value = 500200
$ws.Range("I106").Value = 400
$ws.Range("J106").Value = 1000000
$ws.Range("K106").Value = 400
$ws.Range("L106").Value = 1000000
$ws.Range("M106").Value = 231
$ws.Range("N106").Value = -1001262

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5143.359
$ws.Range("I32").Value = 4620.816
$ws.Range("K32").Value = 4620.816
$ws.Range("M32").Value = -4333.816

$ws.Range("H61").Value = 1219.7632
$ws.Range("I61").Value = 1144.6216
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 1144.6216
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -932.6215999999999
$ws.Range("N61").Value = -4424

$ws.Range("H74").Value = 2894.52
$ws.Range("I74").Value = 2906.5366
$ws.Range("J74").Value = 2839.7778
$ws.Range("K74").Value = 2906.5366
$ws.Range("L74").Value = 2839.7778
$ws.Range("M74").Value = -2032.5366
$ws.Range("N74").Value = -4587.7778

$ws.Range("H77").Value = 2894.52
$ws.Range("I77").Value = 2906.5366
$ws.Range("J77").Value = 2839.7778
$ws.Range("K77").Value = 14532.683
$ws.Range("L77").Value = 14198.889
$ws.Range("M77").Value = -10164.683
$ws.Range("N77").Value = -22934.889

$ws.Range("H122").Value = 2609.1765
$ws.Range("I122").Value = 2439.2727
$ws.Range("K122").Value = 7317.8181
$ws.Range("M122").Value = -4867.8181

$ws.Range("H132").Value = 1834.8214
$ws.Range("I132").Value = 1828.7037
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 5486.1111
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -2956.1111
$ws.Range("N132").Value = -11060

$ws.Range("H136").Value = 1219.7632
$ws.Range("I136").Value = 1144.6216
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 3433.8648
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -883.8647999999998
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()

$ws.Range("H96").Value = 24280.6
$ws.Range("I96").Value = 24280.6
$ws.Range("K96").Value = 24280.6
$ws.Range("M96").Value = -21534.6

$ws.Range("H107").Value = 1609.2059
$ws.Range("I107").Value = 1433.7667
$ws.Range("J107").Value = 2925
$ws.Range("K107").Value = 1433.7667
$ws.Range("L107").Value = 2925
$ws.Range("M107").Value = 486.2333000000001
$ws.Range("N107").Value = -6765

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6534.9487
$ws.Range("I31").Value = 121449
$ws.Range("J31").Value = 3510.8948
$ws.Range("K31").Value = 121449
$ws.Range("L31").Value = 3510.8948
$ws.Range("M31").Value = -121154
$ws.Range("N31").Value = -4100.8948

$ws.Range("H34").Value = 6534.9487
$ws.Range("I34").Value = 121449
$ws.Range("J34").Value = 3510.8948
$ws.Range("K34").Value = 121449
$ws.Range("L34").Value = 3510.8948
$ws.Range("M34").Value = -121247
$ws.Range("N34").Value = -3914.8948

$ws.Range("H58").Value = 1688.6364
$ws.Range("J58").Value = 5266.3335
$ws.Range("L58").Value = 5266.3335
$ws.Range("N58").Value = -5672.3335

$ws.Range("H119").Value = 27500
$ws.Range("J119").Value = 27500
$ws.Range("L119").Value = 27500
$ws.Range("N119").Value = -37176

$ws.Range("H122").Value = 1850.2759
$ws.Range("I122").Value = 2506.353
$ws.Range("K122").Value = 7519.059
$ws.Range("M122").Value = -5069.059

$ws.Range("H132").Value = 2218.081
$ws.Range("I132").Value = 2144.8572
$ws.Range("J132").Value = 3499.5
$ws.Range("K132").Value = 6434.571599999999
$ws.Range("L132").Value = 10498.5
$ws.Range("M132").Value = -3904.571599999999
$ws.Range("N132").Value = -15558.5

$ws.Range("H136").Value = 1688.6364
$ws.Range("J136").Value = 5266.3335
$ws.Range("L136").Value = 15799.0005
$ws.Range("N136").Value = -20899.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1751.4667
$ws.Range("I8").Value = 1751.4667
$ws.Range("K8").Value = 5254.4001
$ws.Range("M8").Value = -5115.4001

$ws.Range("H10").Value = 99.38461
$ws.Range("I10").Value = 99.38461
$ws.Range("K10").Value = 298.15383
$ws.Range("M10").Value = -159.15383

$ws.Range("H12").Value = 245.5
$ws.Range("I12").Value = 300.4
$ws.Range("J12").Value = 215
$ws.Range("K12").Value = 901.1999999999999
$ws.Range("L12").Value = 645
$ws.Range("M12").Value = -728.1999999999999
$ws.Range("N12").Value = -991

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1814.2069
$ws.Range("I102").Value = 1619.7307
$ws.Range("K102").Value = 1619.7307
$ws.Range("M102").Value = 2.26929999999993

$ws.Range("H122").Value = 2704
$ws.Range("I122").Value = 2328
$ws.Range("K122").Value = 6984
$ws.Range("M122").Value = -4534

$ws.Range("H132").Value = 3370.7463
$ws.Range("I132").Value = 2828.1355
$ws.Range("J132").Value = 7372.5
$ws.Range("K132").Value = 8484.406499999999
$ws.Range("L132").Value = 22117.5
$ws.Range("M132").Value = -5954.406499999999
$ws.Range("N132").Value = -27177.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 902.8837
$ws.Range("I16").Value = 963
$ws.Range("K16").Value = 963
$ws.Range("M16").Value = -793

$ws.Range("H22").Value = 4614024.5
$ws.Range("I22").Value = 1570
$ws.Range("J22").Value = 16145161
$ws.Range("K22").Value = 1570
$ws.Range("L22").Value = 16145161
$ws.Range("M22").Value = -1275
$ws.Range("N22").Value = -16145751

$ws.Range("H27").Value = 4614024.5
$ws.Range("I27").Value = 1570
$ws.Range("J27").Value = 16145161
$ws.Range("K27").Value = 1570
$ws.Range("L27").Value = 16145161
$ws.Range("M27").Value = -1463
$ws.Range("N27").Value = -16145375

$ws.Range("H100").Value = 924.5
$ws.Range("I100").Value = 924.5
$ws.Range("K100").Value = 924.5
$ws.Range("M100").Value = -383.5

$ws.Range("H103").Value = 50000
$ws.Range("J103").Value = 50000
$ws.Range("L103").Value = 50000
$ws.Range("N103").Value = -52344

$ws.Range("H132").Value = 4763.222
$ws.Range("I132").Value = 2329.6667
$ws.Range("J132").Value = 8170.2
$ws.Range("K132").Value = 6989.000100000001
$ws.Range("L132").Value = 24510.6
$ws.Range("M132").Value = -4459.000100000001
$ws.Range("N132").Value = -29570.6

$ws.Range("H136").Value = 2158.5588
$ws.Range("J136").Value = 4557.4
$ws.Range("L136").Value = 13672.2
$ws.Range("N136").Value = -18772.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 9096785
$ws.Range("J81").Value = 16675335
$ws.Range("L81").Value = 33350670
$ws.Range("N81").Value = -33352792

$ws.Range("H84").Value = 9096785
$ws.Range("J84").Value = 16675335
$ws.Range("L84").Value = 166753350
$ws.Range("N84").Value = -166763958

$ws.Range("H112").Value = 35587.5
$ws.Range("J112").Value = 35587.5
$ws.Range("L112").Value = 35587.5
$ws.Range("N112").Value = -38541.5

$ws.Range("H119").Value = 43473.5
$ws.Range("J119").Value = 40464.668
$ws.Range("L119").Value = 40464.668
$ws.Range("N119").Value = -50140.668

$ws.Range("H122").Value = 2108.158
$ws.Range("I122").Value = 2108.158
$ws.Range("K122").Value = 6324.474
$ws.Range("M122").Value = -3874.474

$ws.Range("H136").Value = 4443.6553
$ws.Range("J136").Value = 6556.0835
$ws.Range("L136").Value = 19668.2505
$ws.Range("N136").Value = -24768.2505
